# Auto-generated edit script: add 2022-Q1 sheet, rebuild 总计 sheet
$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing "总计" sheet to "2022-Q1" (keeps its sheetId / physical
#     part, matching the diff where sheetId=6 becomes the "2022-Q1" tab), then add a
#     brand-new sheet named "总计" positioned right after it (gets the next sheetId=7). ---
$zongji = $wb.Worksheets.Item("总计")
$zongji.Name = "2022-Q1"
$newZongji = $wb.Worksheets.Add()
$newZongji.Name = "总计"
$q1sheet = $wb.Worksheets.Item("2022-Q1")
$newZongji.Move($null, $q1sheet)

$q1 = $wb.Worksheets.Item("2022-Q1")
$tot = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

# --- Step 2: the "2022-Q1" sheet currently still holds the old 4-column 总计 table;
#     clear it completely before rebuilding it as an 8-column fund-holdings detail sheet,
#     the same shape used by the other quarterly detail tabs (e.g. "2021-Q4"). ---
$q1.Range("A1:H40").Clear()

# --- Step 3: clone header / index-column cell formatting from the "2021-Q4" template so the
#     new sheet keeps the same bold-centered-bordered style (cellXfs index 2) instead of
#     Excel inventing a brand-new style entry. ---
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:A35").Copy()
$q1.Range("A2:A35").PasteSpecial(-4122)

# --- Step 4: header row text (B1:H1) ---
$q1.Range("B1").Formula = '="基金代码"'
$q1.Range("C1").Formula = '="基金名称"'
$q1.Range("D1").Formula = '="基金规模"'
$q1.Range("E1").Formula = '="股票总仓位"'
$q1.Range("F1").Formula = '="仓位占比"'
$q1.Range("G1").Formula = '="持有市值(亿元)"'
$q1.Range("H1").Formula = '="仓位排名"'
$q1.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4163)

# --- Step 5: fund holdings detail rows (2..35) ---
#     columns B..G are text cells in the source workbook (values such as "35.89" are
#     stored as strings, not numbers) so each is written via a ="..." formula and then
#     flattened to a static value with PasteSpecial(xlPasteValues); A (0-based index) and
#     H (rank) are genuine numbers written directly. ---
$fundData = @(
    @('010454','交银施罗德内需增长一年持有期混合','35.89','91.74','8.09','2.9035','4'),
    @('012582','交银施罗德品质增长一年持有期混合型证券投资基金A','51.62','92.40','5.48','2.8288','4'),
    @('001216','易方达新收益灵活配置混合 - A','68.57','91.39','3.09','2.1188','10'),
    @('519714','交银施罗德消费新驱动股票','18.40','91.63','8.99','1.6542','4'),
    @('004868','交银施罗德股息优化混合','16.29','92.33','8.93','1.4547','4'),
    @('004075','交银施罗德医药创新股票','31.60','83.63','4.54','1.4346','7'),
    @('005004','交银施罗德品质升级混合','15.90','93.72','9.00','1.4310','6'),
    @('009618','交银施罗德启汇混合','21.00','71.65','4.52','0.9492','3'),
    @('260110','景顺长城精选蓝筹混合','17.93','92.90','5.23','0.9377','7'),
    @('090001','大成价值增长混合','18.65','61.32','4.07','0.7591','6'),
    @('010275','嘉实优质精选混合A','15.40','92.06','4.60','0.7084','8'),
    @('001217','易方达新收益灵活配置混合 - C','19.83','91.39','3.09','0.6127','10'),
    @('011335','银河医药健康混合型证券投资基金','8.94','87.69','6.44','0.5757','6'),
    @('519710','交银施罗德策略回报灵活配置混合','5.77','79.92','8.02','0.4628','6'),
    @('010405','惠升医药健康6个月持有期混合','12.18','64.63','2.90','0.3532','8'),
    @('159883','永赢中证全指医疗器械交易型开放式指数证券投资基金','9.17','99.23','3.05','0.2797','6'),
    @('260112','景顺长城能源基建混合','16.49','60.89','1.45','0.2391','10'),
    @('090020','大成健康产业混合','3.76','91.73','5.78','0.2173','7'),
    @('012045','大成医药健康股票A','2.87','93.58','7.25','0.2081','5'),
    @('519655','银河现代服务主题灵活配置混合','3.89','85.22','4.39','0.1708','8'),
    @('007066','浦银安盛先进制造混合A','3.66','74.55','3.91','0.1431','7'),
    @('519013','海富通风格优势混合','3.57','89.46','3.76','0.1342','6'),
    @('012583','交银施罗德品质增长一年持有期混合型证券投资基金C','2.37','92.40','5.48','0.1299','4'),
    @('007067','浦银安盛先进制造混合C','2.27','74.55','3.91','0.0888','7'),
    @('000587','大成灵活配置混合','2.39','85.02','3.16','0.0755','10'),
    @('159898','招商中证全指医疗器械交易型开放式指数证券投资基金','1.61','99.41','3.03','0.0488','6'),
    @('010276','嘉实优质精选混合C','1.01','92.06','4.60','0.0465','8'),
    @('519139','海富通沪港深灵活配置混合','1.32','94.37','3.15','0.0416','9'),
    @('159873','天弘中证全指医疗保健设备与服务ETF','1.43','99.59','2.64','0.0378','8'),
    @('159891','建信中证全指医疗保健设备与服务交易型开放式指数证券投资基金','1.40','95.24','2.53','0.0354','8'),
    @('013441','西藏东财创新医疗六个月定开混合','0.58','81.46','4.81','0.0279','10'),
    @('012046','大成医药健康股票C','0.25','93.58','7.25','0.0181','5'),
    @('516610','大成中证全指医疗保健设备与服务交易型开放式指数证券投资基金','0.69','96.20','2.49','0.0172','8'),
    @('006992','嘉合锦创优势精选混合','0.02','74.79','4.66','0.0009','3'),
)
for ($i = 0; $i -lt $fundData.Count; $i++) {
    $row = $fundData[$i]
    $r = $i + 2
    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).Formula = '="' + $row[0] + '"'
    $q1.Cells.Item($r, 3).Formula = '="' + $row[1] + '"'
    $q1.Cells.Item($r, 4).Formula = '="' + $row[2] + '"'
    $q1.Cells.Item($r, 5).Formula = '="' + $row[3] + '"'
    $q1.Cells.Item($r, 6).Formula = '="' + $row[4] + '"'
    $q1.Cells.Item($r, 7).Formula = '="' + $row[5] + '"'
    $q1.Cells.Item($r, 8).Value = $row[6]
}
$q1.Range("B2:G35").Copy()
$q1.Range("B2:G35").PasteSpecial(-4163)

# --- Step 6: rebuild the "总计" sheet: header row + one summary row per quarter, with the
#     new "2022-Q1" row inserted at the top (index 0) and the rest pushed down. ---
$q1.Range("A2:A6").Copy()
$tot.Range("A2:A7").PasteSpecial(-4122)
$q1.Range("B1:B1").Copy()
$tot.Range("B1:D1").PasteSpecial(-4122)

$tot.Range("B1").Formula = '="日期"'
$tot.Range("C1").Formula = '="持有数量(只)"'
$tot.Range("D1").Formula = '="持有市值(亿元)"'
$tot.Range("B1:D1").Copy()
$tot.Range("B1:D1").PasteSpecial(-4163)

$summaryData = @(
    @('2022-Q1','34','21.15'),
    @('2021-Q4','80','50.26'),
    @('2021-Q3','63','32.57'),
    @('2021-Q2','146','101.94'),
    @('2021-Q1','184','58.83'),
    @('2020-Q4','197','47.95'),
)
for ($i = 0; $i -lt $summaryData.Count; $i++) {
    $row = $summaryData[$i]
    $r = $i + 2
    $tot.Cells.Item($r, 1).Value = $i
    $tot.Cells.Item($r, 3).Value = $row[1]
    $tot.Cells.Item($r, 4).Value = $row[2]
}
$dateData = @(
    '2022-Q1',
    '2021-Q4',
    '2021-Q3',
    '2021-Q2',
    '2021-Q1',
    '2020-Q4',
)
for ($i = 0; $i -lt $dateData.Count; $i++) {
    $r = $i + 2
    $tot.Cells.Item($r, 2).Formula = '="' + $dateData[$i] + '"'
}
$tot.Range("B2:B7").Copy()
$tot.Range("B2:B7").PasteSpecial(-4163)
